# Fixes GitHub issue #3: cells that carry a type but no value (e.g. an
# empty date/number cell) must still round-trip without raising.
#
# Adds a new "Empty Eagress" row to the Posts sheet whose Title/Comment
# cells are left without values, while the Date/Comment Count cells keep
# their number formatting despite also being empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Posts")

# New header-ish label + explanatory comment (become new shared strings).
$ws.Range("A5").Value = "Empty Eagress"
$ws.Range("C5").Value = "The title, date, and comment have types, but no values"

# B5 is intentionally left blank.
# D5 / E5 stay empty but retain the column's number formats (date / 0.00),
# exercising the "typed cell with no value" case from the bug report.
$ws.Range("D5").NumberFormat = "mm-dd-yy"
$ws.Range("E5").NumberFormat = "0.00"
